$wb = $excel.ActiveWorkbook

# --- Overview sheet: handback status text + widened columns E/F ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527

# --- zh-cn sheet: status text, handback file/link/datetime, wider columns ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("I2").Value = "87f633d4-f3dc-4621-978f-1a8809691b23.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c946e04ce15585c6cce410a98c74454e5eea5f0e/e2e/87f633d4-f3dc-4621-978f-1a8809691b23.md", "", "", "87f633d4-f3dc-4621-978f-1a8809691b23.md")
$wsZhCn.Range("J2").Value = "87f633d4-f3dc-4621-978f-1a8809691b23.a8545c68db732f3a4a7574ee2210eeae79ff480c.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-28 02:58:44"
$wsZhCn.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsZhCn.Columns.Item(9).ColumnWidth = 40
$wsZhCn.Columns.Item(10).ColumnWidth = 40

# --- de-de sheet: status text, handback file/link/datetime, wider columns ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("I2").Value = "87f633d4-f3dc-4621-978f-1a8809691b23.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c946e04ce15585c6cce410a98c74454e5eea5f0e/e2e/87f633d4-f3dc-4621-978f-1a8809691b23.md", "", "", "87f633d4-f3dc-4621-978f-1a8809691b23.md")
$wsDeDe.Range("J2").Value = "87f633d4-f3dc-4621-978f-1a8809691b23.a8545c68db732f3a4a7574ee2210eeae79ff480c.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-28 02:58:51"
$wsDeDe.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsDeDe.Columns.Item(9).ColumnWidth = 40
$wsDeDe.Columns.Item(10).ColumnWidth = 40
